# Clean up the item list: drop the extra rows (old rows 15-21) and fix the
# duplicated "tag" value / price typo that had crept into row 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix price typo and the tag-column off-by-one (rows 7-14 each re-point to
# the previous "waterN" tag, i.e. row7 -> water5, row8 -> water6, ... row14 -> water12)
$ws.Range("C7").Value = 4000
$ws.Range("D7").Value = "water5"
$ws.Range("D8").Value = "water6"
$ws.Range("D9").Value = "water7"
$ws.Range("D10").Value = "water8"
$ws.Range("D11").Value = "water9"
$ws.Range("D12").Value = "water10"
$ws.Range("D13").Value = "water11"
$ws.Range("D14").Value = "water12"

# Remove the now-unwanted trailing rows (old rows 15-21)
$ws.Rows("15:21").Delete()

# Match the author's final selection
$null = $ws.Range("B15").Select()
